$d = $word.ActiveDocument

function Set-RangeXml($rng, $innerXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Part 1: the "If we count even further to 20" list currently uses numId 19;
# it must be switched to numId 13 (the same list definition used by the
# "count to 10" list above it), one paragraph at a time so the rest of the
# document is left untouched.
# ---------------------------------------------------------------------------
$words = @("Middle", "Ring", "Pinky", "Ring", "Middle", "Index", "Thumb", "Index", "Middle")

# Locate the first paragraph of that second list by scanning for the
# paragraph right after "If we count even further to 20, the pattern goes as such:"
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "If we count even further to 20") {
        $anchorIdx = $i + 1
        break
    }
}
if ($anchorIdx -lt 0) {
    throw "Could not locate the 'count even further to 20' list anchor paragraph"
}

for ($k = 0; $k -lt $words.Count; $k++) {
    $p = $d.Paragraphs.Item($anchorIdx + $k)
    $w = $words[$k]
    $inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> ' + $w + '</w:t></w:r></w:p>'
    Set-RangeXml $p.Range $inner
}

# ---------------------------------------------------------------------------
# Part 2: the final two paragraphs of the document -- the last "Ring" list
# item (numId 19 -> 13, same as above) and the trailing empty NormalWeb
# paragraph that carries the _GoBack bookmark -- are replaced by:
#   1) the last "Ring" list item (now numId 13)
#   2) two blank paragraphs
#   3) a new paragraph (carrying the _GoBack bookmark) with the "Notice
#      that..." commentary, using italics for the finger names
#   4) a new trailing blank NormalWeb paragraph
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$pPenultimate = $d.Paragraphs.Item($n - 1)
$pLast = $d.Paragraphs.Item($n)
$tailRange = $d.Range($pPenultimate.Range.Start, $pLast.Range.End)

$dash = [char]0x2014
$finalInner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> Ring</w:t></w:r></w:p>' +
    '<w:p/>' +
    '<w:p/>' +
    '<w:p>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t xml:space="preserve">Notice that the </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>Ring</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and the </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Index </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">finger are counted </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>twice</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> in one round' + $dash + 'one round counting every 10 to 5 fingers: therefore, there is a 2 of 10 chance of every tenth number to land on either the </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Ring </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">or </w:t></w:r>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">Index </w:t></w:r>' +
    '<w:r><w:t>finger.</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr></w:p>'

Set-RangeXml $tailRange $finalInner

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
